$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2048
$ws1.Range("F6").Value = 3508
$ws1.Range("F9").Value = 1458
$ws1.Range("F10").Value = 4856
$ws1.Range("F12").Value = 1810
$ws1.Range("F22").Value = 100
$ws1.Range("F27").Value = 1167
$ws1.Range("F31").Value = 508
$ws1.Range("F34").Value = 1868
$ws1.Range("F36").Value = 1101
$ws1.Range("F42").Value = 540
$ws1.Range("F47").Value = 530

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 128

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2048
$ws4.Range("F6").Value = 3508
$ws4.Range("F8").Value = 1458
$ws4.Range("F9").Value = 4856
$ws4.Range("F10").Value = 1810
$ws4.Range("F22").Value = 100
$ws4.Range("F26").Value = 1167
$ws4.Range("F32").Value = 1868
$ws4.Range("F34").Value = 1101
$ws4.Range("F42").Value = 540
$ws4.Range("F45").Value = 530
